$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update "Dato lagret" value (force text so Excel doesn't turn it into a date serial,
# then clear the format so the cell keeps its original unstyled appearance)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2021-02-04"
$ws.Range("B3").ClearFormats()

# Row 4 used to be "egenskap" / formula text, now becomes "overlapp" / new formula text
$ws.Range("A4").Value = "overlapp"
$ws.Range("B4").Value = "591(5277 < 4 AND (5270=8168 OR 5270=8149))"

# Row 5 used to be "overlapp" / "60(1263=7304)", now becomes "egenskap" / "1263=7304"
$ws.Range("A5").Value = "egenskap"
$ws.Range("B5").Value = "1263=7304"

# Remove rows 8 and 9 (egenskapfilter_bru and the overlapp search row), shrinking the sheet to A1:B7
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
